# Applies the "Rejestracja" (registration) test-sheet edit:
#  - Fixes the ID ordering of the existing rows 3 & 4 (swap 3<->2)
#  - Inserts a brand-new test-case row ("existing username + existing mail + password")
#    as the new row 6, pushing the remaining rows down
#  - Renumbers the ID column for the shifted rows
#  - Extends the "Tabela32" table / autofilter to the new range
#  - Widens column B and updates the sheet view (top-left cell & selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Rejestracja"

# --- Fix ID values on existing rows 3 and 4 (they were swapped) -----------
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# --- Insert a new row at position 6, shifting the old rows 6-8 down -------
$ws.Rows.Item(6).Insert()

# Copy formatting (cell style) for the three new data cells from the row
# that is now directly below (row 7, which used to be row 6 and carries the
# correct "existing username/mail" style s="2").
$ws.Range("A7").Copy($ws.Range("A6"))
$ws.Range("B7").Copy($ws.Range("B6"))
$ws.Range("C7").Copy($ws.Range("C6"))

# --- Populate the new row's content ----------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Rejestracja istniejący użytkownik+ istniejący mail+ hasło"
$ws.Range("C6").Value = "Sprawdza czy wyskakuje:`n2 errors`n    Username has already been taken`n    Email has already been taken"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "Damian Lechański"

$ws.Rows.Item(6).RowHeight = 75

# --- Renumber the ID column for the rows that shifted down -----------------
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6

$ws.Rows.Item(8).RowHeight = 30

# --- Extend the table / autofilter range to include the new row ------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E9"))

# --- Column B width ----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 53.33

# --- Sheet view: selection --------------------------------------------------
$ws.Range("E5:E6").Select()
